# Edit TestData.xlsx: Employee_Details sheet
# - Insert a new column before BK (alt work location address line 2)
# - Add new "emergency contact" columns BO:BW with header row + sample data row 8
# - Add mailto hyperlink on the emergency contact email cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employee_Details")

# 1) Insert a new column at BK; this shifts the existing
#    altWorkLocationAddressType/Country/ZipCode columns (and their row7 data)
#    one column to the right (BK:BM -> BL:BN).
$ws.Columns("BK").Insert()

# 2) New header labels for the emergency contact columns (written in the
#    same order the authoring session introduced them).
$ws.Range("BP1").Value2 = "familyName"
$ws.Range("BQ1").Value2 = "relationship"
$ws.Range("BS1").Value2 = "phoneType"
$ws.Range("BU1").Value2 = "phoneNumber"
$ws.Range("BV1").Value2 = "emailType"

# 3) New "alternative work location" second address line.
$ws.Range("BK7").Value2 = "Adamo Ave"
$ws.Range("BK1").Value2 = "altWorkLocationAddressLine2"

# 4) Emergency-contact header cells that reuse already-existing shared strings.
$ws.Range("BO1").Value2 = "firstName"
$ws.Range("BR1").Value2 = "gender"
$ws.Range("BT1").Value2 = "countryCode"
$ws.Range("BW1").Value2 = "email"

# 5) Sample emergency-contact data row (row 8).
$ws.Range("BO8").Value2 = "Thomas"
$ws.Range("BP8").Value2 = "Edward"
$ws.Range("BQ8").Value2 = "Brother"
$ws.Range("BR8").Value2 = "Male"
$ws.Range("BS8").Value2 = "Work Mobile Phone"
$ws.Range("BT8").Value2 = "United States 1"
$ws.Range("BV8").Value2 = "Work Email"

# Email cell carries a mailto hyperlink (gets the built-in Hyperlink style).
$ws.Range("BW8").Value2 = "Thomas.Edward@cognizant.com"
$ws.Hyperlinks.Add($ws.Range("BW8"), "mailto:Thomas.Edward@cognizant.com")

# Phone number kept as text with a leading apostrophe (quote-prefixed string).
$ws.Range("BU8").Value2 = "'6323145"

# 6) Match the saved scroll/selection state.
$ws.Activate()
$ws.Range("BQ5").Select()
